# Add two new date columns (19-sep, 20-sep) to the right of the existing
# data table, with header labels in row 1 and numeric values for each
# product row (2-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new date labels.
# Copy the existing last header cell (BS1) first so the new header cells
# inherit the same text style used throughout row 1.
$ws.Range("BS1").Copy()
$ws.Range("BT1:BU1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("BT1").Value = "19-sep"
$ws.Range("BU1").Value = "20-sep"

# Copy the existing last data cell (BS2) formatting across the new data
# range so the new columns match the integer/centered style used by the
# rest of the table.
$ws.Range("BS2").Copy()
$ws.Range("BT2:BU11").PasteSpecial(-4122)  # xlPasteFormats

# Data rows 2-11 for the two new columns
$values = @{
    2  = @(13, 14)
    3  = @(11, 11)
    4  = @(12, 12)
    5  = @(10, 10)
    6  = @(12, 11)
    7  = @(17, 16)
    8  = @(18, 17)
    9  = @(12, 11)
    10 = @(18, 18)
    11 = @(6, 5)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 72).Value = $pair[0]  # column BT
    $ws.Cells.Item($row, 73).Value = $pair[1]  # column BU
}

# Update the active selection to reflect the post-edit state
$ws.Range("BZ12").Select()
